$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name corrections in column A (re-sorted ranking) ---
$ws.Range("A1").Value = 'Datos actualizados a 6 de Mayo de 2020 a las 15:03'
$ws.Range("A20").Value = 'Arabia Saudita'
$ws.Range("A21").Value = 'Ecuador'
$ws.Range("A23").Value = 'Portugal'
$ws.Range("A24").Value = 'Mexico'
$ws.Range("A61").Value = 'Barein'
$ws.Range("A62").Value = 'Luxemburgo'
$ws.Range("A72").Value = 'Camerun'
$ws.Range("A73").Value = 'Uzbekistan'
$ws.Range("A74").Value = 'Croacia'
$ws.Range("A122").Value = 'Guinea Ecuatorial'
$ws.Range("A123").Value = 'Taiwan'
$ws.Range("A124").Value = 'Paraguay'
$ws.Range("A125").Value = 'Reunion'
$ws.Range("A126").Value = 'Gabon'
$ws.Range("A127").Value = 'Estado de Palestina'
$ws.Range("A128").Value = 'Venezuela'
$ws.Range("A129").Value = 'Mauricio'
$ws.Range("A130").Value = 'Isla de Man'
$ws.Range("A131").Value = 'Montenegro'
$ws.Range("A145").Value = 'Madagascar'
$ws.Range("A146").Value = 'Guadalupe'
$ws.Range("A162").Value = 'Republica de Africa Central'
$ws.Range("A163").Value = 'Guyana'
$ws.Range("A164").Value = 'Bahamas'
$ws.Range("A173").Value = 'Siria'
$ws.Range("A174").Value = 'Macao'
$ws.Range("A205").Value = 'Montserrat'
$ws.Range("A206").Value = 'Seychelles'

# --- Updated statistics (columns B-H) ---
$ws.Range("B4").Value = 1238463
$ws.Range("C4").Value = 830
$ws.Range("E4").Value = 965165
$ws.Range("G4").Value = 16
$ws.Range("H4").Value = 72287

$ws.Range("B20").Value = 31938
$ws.Range("C20").Value = 1687
$ws.Range("D20").Value = 6783
$ws.Range("E20").Value = 24946
$ws.Range("F20").Value = 143
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 209

$ws.Range("B21").Value = 31881
$ws.Range("D21").Value = 3433
$ws.Range("E21").Value = 26879
$ws.Range("F21").Value = 159
$ws.Range("H21").Value = 1569

$ws.Range("B22").Value = 30060
$ws.Range("C22").Value = 51
$ws.Range("E22").Value = 2865

$ws.Range("B23").Value = 26182
$ws.Range("C23").Value = 480
$ws.Range("D23").Value = 2076
$ws.Range("E23").Value = 23017
$ws.Range("F23").Value = 136
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = 1089

$ws.Range("B24").Value = 26025
$ws.Range("C24").Value = 1120
$ws.Range("D24").Value = 16810
$ws.Range("E24").Value = 6708
$ws.Range("F24").Value = 378
$ws.Range("G24").Value = 236
$ws.Range("H24").Value = 2507

$ws.Range("B25").Value = 23918
$ws.Range("C25").Value = 702
$ws.Range("E25").Value = 16903
$ws.Range("F25").Value = 425
$ws.Range("G25").Value = 87
$ws.Range("H25").Value = 2941

$ws.Range("D43").Value = 7493
$ws.Range("E43").Value = 1939
$ws.Range("F43").Value = 46
$ws.Range("G43").Value = 3
$ws.Range("H43").Value = 506

$ws.Range("B54").Value = 6289
$ws.Range("C54").Value = 485
$ws.Range("D54").Value = 2219
$ws.Range("E54").Value = 4028
$ws.Range("F54").Value = 88
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 42

$ws.Range("B60").Value = 4344
$ws.Range("C60").Value = 139
$ws.Range("E60").Value = 2927

$ws.Range("B61").Value = 3842
$ws.Range("C61").Value = 122
$ws.Range("D61").Value = 1860
$ws.Range("E61").Value = 1974
$ws.Range("F61").Value = 4
$ws.Range("H61").Value = 8

$ws.Range("B62").Value = 3840
$ws.Range("D62").Value = 3412
$ws.Range("E62").Value = 332
$ws.Range("F62").Value = 22
$ws.Range("H62").Value = 96

$ws.Range("B72").Value = 2265
$ws.Range("C72").Value = 161
$ws.Range("D72").Value = 1000
$ws.Range("E72").Value = 1201
$ws.Range("F72").Value = 12
$ws.Range("H72").Value = 64

$ws.Range("B73").Value = 2231
$ws.Range("C73").Value = 24
$ws.Range("D73").Value = 1556
$ws.Range("E73").Value = 665
$ws.Range("F73").Value = 8
$ws.Range("H73").Value = 10

$ws.Range("B74").Value = 2119
$ws.Range("C74").Value = 7
$ws.Range("D74").Value = 1601
$ws.Range("E74").Value = 433
$ws.Range("F74").Value = 14
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 85

$ws.Range("B101").Value = 774
$ws.Range("C101").Value = 3
$ws.Range("E101").Value = 550

$ws.Range("B115").Value = 574
$ws.Range("C115").Value = 1
$ws.Range("E115").Value = 552

$ws.Range("C122").Value = 124
$ws.Range("D122").Value = 13
$ws.Range("E122").Value = 422
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 4

$ws.Range("B123").Value = 439
$ws.Range("C123").Value = 1
$ws.Range("D123").Value = 339
$ws.Range("E123").Value = 94
$ws.Range("F123").Value = 0
$ws.Range("H123").Value = 6

$ws.Range("B124").Value = 431
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 135
$ws.Range("E124").Value = 286
$ws.Range("F124").Value = 7
$ws.Range("H124").Value = 10

$ws.Range("B125").Value = 425
$ws.Range("C125").Value = 1
$ws.Range("D125").Value = 300
$ws.Range("E125").Value = 125
$ws.Range("F125").Value = 2
$ws.Range("H125").Value = 0

$ws.Range("B126").Value = 397
$ws.Range("D126").Value = 93
$ws.Range("E126").Value = 298
$ws.Range("F126").Value = 1
$ws.Range("H126").Value = 6

$ws.Range("B127").Value = 371
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 127
$ws.Range("E127").Value = 242
$ws.Range("F127").Value = 0
$ws.Range("H127").Value = 2

$ws.Range("B128").Value = 367
$ws.Range("C128").Value = 6
$ws.Range("D128").Value = 164
$ws.Range("E128").Value = 193
$ws.Range("F128").Value = 1

$ws.Range("B129").Value = 332
$ws.Range("D129").Value = 319
$ws.Range("E129").Value = 3
$ws.Range("F129").Value = 3
$ws.Range("H129").Value = 10

$ws.Range("B130").Value = 326
$ws.Range("D130").Value = 271
$ws.Range("E130").Value = 32
$ws.Range("F130").Value = 19
$ws.Range("H130").Value = 23

$ws.Range("B131").Value = 324
$ws.Range("D131").Value = 261
$ws.Range("E131").Value = 55
$ws.Range("F131").Value = 2
$ws.Range("H131").Value = 8

$ws.Range("B136").Value = 225
$ws.Range("C136").Value = 26
$ws.Range("D136").Value = 54
$ws.Range("E136").Value = 157
$ws.Range("G136").Value = 3
$ws.Range("H136").Value = 14

$ws.Range("B145").Value = 158
$ws.Range("C145").Value = 7
$ws.Range("D145").Value = 101
$ws.Range("E145").Value = 57
$ws.Range("F145").Value = 1
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 0

$ws.Range("B146").Value = 152
$ws.Range("D146").Value = 104
$ws.Range("E146").Value = 35
$ws.Range("F146").Value = 4
$ws.Range("G146").Value = 1
$ws.Range("H146").Value = 13

$ws.Range("D151").Value = 77
$ws.Range("E151").Value = 42

$ws.Range("B162").Value = 94
$ws.Range("C162").Value = 9
$ws.Range("D162").Value = 10
$ws.Range("E162").Value = 84
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 0

$ws.Range("B163").Value = 93
$ws.Range("C163").Value = 1
$ws.Range("D163").Value = 27
$ws.Range("E163").Value = 56
$ws.Range("F163").Value = 3
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 10

$ws.Range("B164").Value = 89
$ws.Range("D164").Value = 26
$ws.Range("E164").Value = 52
$ws.Range("F164").Value = 1
$ws.Range("H164").Value = 11

$ws.Range("C173").Value = 1
$ws.Range("D173").Value = 27
$ws.Range("E173").Value = 15
$ws.Range("F173").Value = 0
$ws.Range("H173").Value = 3

$ws.Range("B174").Value = 45
$ws.Range("D174").Value = 39
$ws.Range("E174").Value = 6
$ws.Range("F174").Value = 1
$ws.Range("H174").Value = 0

$ws.Range("D205").Value = 7
$ws.Range("F205").Value = 1
$ws.Range("H205").Value = 1

$ws.Range("D206").Value = 8
$ws.Range("F206").Value = 0
$ws.Range("H206").Value = 0

